# Fixed typo in paper and edited figure 1
#
# 1) Update the "datetimeFigureOut" date placeholder (2/23/2023 -> 5/12/2023)
#    on the slide master and every slide layout.
# 2) Update the Figure 1 title text box on slide 1.

$p = $ppt.ActivePresentation

$oldDate = "2/23/2023"
$newDate = "5/12/2023"

# --- Slide master date placeholder ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Every slide layout's date placeholder ---
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Figure 1 title textbox on slide 1 ---
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $tr = $sh.TextFrame.TextRange
        $firstPara = $tr.Paragraphs(1, 1)
        $firstParaText = $firstPara.Text.TrimEnd("`r")
        if ($firstParaText -eq "Monthly Temperature and Precipitation at 0.5-deg") {
            $firstPara.Text = "Gridded Temperature and Precipitation"
        }
    }
}
